$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns J and K (rows 1-51) were re-sourced: previously J held the
# string "r" (row1) / numeric 1 (rows 2-51) and K held string "s" (row1)
# / numeric 0.5 (rows 2-51). The corrected data source now provides a
# constant 0.5 / 0.3 for every row including the header row, so the
# shared strings are no longer used at all.
$ws.Range("J1:J51").Value = 0.5
$ws.Range("K1:K51").Value = 0.3

# The author also scrolled/reselected before saving.
$excel.ActiveWindow.ScrollRow = 36
$ws.Range("K1:K51").Select()
